$d = $word.ActiveDocument

function Replace-AllText($findText, $replaceText) {
    $range = $d.Content
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2)
}

Replace-AllText "Имя регистра" "Регистр"
Replace-AllText "Имя поля" "Поле"
Replace-AllText "Reset" "Значение"
Replace-AllText "Имя enum" "Enum"
